$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.262.58"
$ws.Range("E2").Value = "  -1.29%  "

$ws.Range("D3").Value = "2.277.32"
$ws.Range("E3").Value = "  -1.67%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "264.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.93%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.91%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("E9").Value = "  -2.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0930"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.03%  "

$ws.Range("E13").Value = "  +1.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.76%  "

$ws.Range("D15").Value = "2.618.71"
$ws.Range("E15").Value = "  -1.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.855"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.87%  "

$ws.Range("D17").Value = "2.273.19"
$ws.Range("E17").Value = "  -1.41%  "

$ws.Range("D18").Value = "43.198.00"
$ws.Range("E18").Value = "  -1.52%  "

$ws.Range("E19").Value = "  -2.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.67"
$ws.Range("D24").Style = "Normal"

$ws.Range("E25").Value = "  -0.99%  "

$ws.Range("E26").Value = "  +0.34%  "

$ws.Range("E27").Value = "  -1.67%  "

$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.02%  "

$ws.Range("E30").Value = "  -3.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0902"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.93%  "

$ws.Range("E36").Value = "  +0.44%  "

$ws.Range("E37").Value = "  -2.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0352"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.22%  "

$ws.Range("E40").Value = "  -6.81%  "

$ws.Range("E41").Value = "  +8.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "76.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.51%  "

$ws.Range("E44").Value = "  -3.26%  "

$ws.Range("E45").Value = "  +1.04%  "

$ws.Range("E46").Value = "  -0.12%  "

$ws.Range("E47").Value = "  -2.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.64%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0993"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.45%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.98%  "
